$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E) for every worker row changes from 2507 to 2508.
# Update all affected rows together so they keep resolving to a single shared value.
$ws.Range("E16:E23").Value = "2508"

# Two workers are removed from the account-statement table:
#   row 19 - MATEO DE JESUS MENDOZA GOMEZ
#   row 20 - JESSE DE JESUS OSORIO CASTELLON
# Deleting the rows shifts the remaining workers (MICHELL, LINDA, MARIA) up.
$ws.Rows("19:20").Delete()

# Update the header summary figures.
$ws.Range("E11").Value = 341640   # VALOR MORA total
$ws.Range("C13").Value = 6        # Cant. Trabajadores

# MARIA FERNANDA ESPINOSA PADILLA's row (now row 21) gets an updated "Valor Mora".
$ws.Range("F21").Value = 56940
